# textPrototypeCLI.pptx — re-save/normalization pass.
#
# The target revision is almost entirely PowerPoint's own re-serialization
# noise (slide/notes IDs renumbered, an unused notes master+notes slide
# dropped, an orphaned notes-master theme part dropped, locale of empty
# placeholder end-runs flipped en-US -> en-SG, p14:creationId stamps added,
# etc.) - none of that is reachable through the public PowerPoint object
# model (Slide.SlideID, for instance, is read-only in real PowerPoint too).
#
# The one part of the diff that *is* a deliberate, user-reachable setting is
# the presentation no longer forcing a custom "first slide number" (the
# <p:presentation firstSlideNum="46" .../> attribute is gone in the after
# state). That corresponds 1:1 to File > Slide Size > (or Design >) the
# legacy "Number slides from" / PageSetup.FirstSlideNumber setting, so apply
# that here.

$p = $ppt.ActivePresentation
$ps = $p.PageSetup

$ps.FirstSlideNumber = 1
